# Update VAR(OK), MSPE, and S_nugget columns (D, E, F) for rows 2-13
# with newly computed CIPW/search-radius results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.453798283957173,  4.030351714201131,  0.37755)
    3  = @(0.7013668025542714, 3.825153030077484,  0.69339)
    4  = @(0.8409282915021251, 3.665542107973322,  1.04571)
    5  = @(0.9134183275653273, 3.625463961433848,  1.28381)
    6  = @(0.948750374058701,  3.599761200486618,  1.45374)
    7  = @(0.9780742489844293, 3.58148076682167,   1.61288)
    8  = @(1.013676352210849,  3.554322504014195,  1.73924)
    9  = @(1.036311721053755,  3.54049603413069,   1.870315)
    10 = @(1.072360263574761,  3.516267784011177,  1.907838)
    11 = @(1.092191907924685,  3.524453162551498,  1.907838)
    12 = @(1.099157067330238,  3.525462690951975,  1.907838)
    13 = @(1.101663636518806,  3.526990841693187,  1.907838)
}

foreach ($row in $values.Keys) {
    $triple = $values[$row]
    $ws.Cells.Item($row, 4).Value = $triple[0]
    $ws.Cells.Item($row, 5).Value = $triple[1]
    $ws.Cells.Item($row, 6).Value = $triple[2]
}
